# "memperbaiki sisa yang masih error (4)" - fix remaining data-type errors
# on the "Flow Database" flow-chart sheet.
#
# Concrete edits:
#   1. F14 ("class" table's 2nd column type) was "varchar", should be "int".
#   2. The "schedule" table block (rows 20-21) had its two attribute rows in
#      the wrong order - swap E20:F20 with E21:F21 so student_identity_id/int
#      comes before class_name/varchar.
#   3. Leave the active selection on E9 (matches the saved cursor position).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. fix the wrong data type on row 14 ------------------------------
$ws.Range("F14").Value = "int"

# --- 2. swap the two attribute rows of the "schedule" table ------------
# Use an out-of-the-way scratch cell as a temp holder so both value and
# cell formatting move together (Range.Copy(Destination) carries format).
$scratch = $ws.Range("Z1:AA1")

$ws.Range("E20:F20").Copy($scratch)
$ws.Range("E21:F21").Copy($ws.Range("E20:F20"))
$scratch.Copy($ws.Range("E21:F21"))
$scratch.Clear()
$excel.CutCopyMode = $false

# --- 3. restore the saved selection -------------------------------------
[void]$ws.Activate()
[void]$ws.Range("E9").Select()
